$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing row 53 values ---
$ws.Range("K53").Value = 82478
$ws.Range("M53").Value = 4584
$ws.Range("N53").Value = 25085
$ws.Range("P53").Value = 22850
$ws.Range("S53").Value = 7043
$ws.Range("T53").Value = 9032
$ws.Range("AF53").Value = 58069
$ws.Range("AJ53").Value = 12586
$ws.Range("AN53").Value = 31249
$ws.Range("BI53").Value = 240426

# --- Add new row 54 ---
# Force the date-like label to be stored as text (matching the style of
# the other "Serie" labels in column A), instead of letting Excel's
# autoconvert turn it into a date serial number.
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = "01-04-2021"
$ws.Range("A54").ClearFormats()
$ws.Range("B54").Value = 31868
$ws.Range("C54").Value = 17846
$ws.Range("D54").Value = 8131
$ws.Range("E54").Value = 142
$ws.Range("F54").Value = 3269
$ws.Range("G54").Value = 128
$ws.Range("H54").Value = 85
$ws.Range("I54").Value = 2257
$ws.Range("J54").Value = 9
$ws.Range("K54").Value = 79625
$ws.Range("L54").Value = 8679
$ws.Range("M54").Value = 4119
$ws.Range("N54").Value = 23726
$ws.Range("O54").Value = 448
$ws.Range("P54").Value = 20865
$ws.Range("Q54").Value = 6447
$ws.Range("R54").Value = 113
$ws.Range("S54").Value = 6287
$ws.Range("T54").Value = 8941
$ws.Range("U54").Value = 0
$ws.Range("V54").Value = 10499
$ws.Range("W54").Value = 447
$ws.Range("X54").Value = 1862
$ws.Range("Y54").Value = 124
$ws.Range("Z54").Value = 1424
$ws.Range("AA54").Value = 2868
$ws.Range("AB54").Value = 0
$ws.Range("AC54").Value = 720
$ws.Range("AD54").Value = 3056
$ws.Range("AE54").Value = 27
$ws.Range("AF54").Value = 58632
$ws.Range("AG54").Value = 401
$ws.Range("AH54").Value = 2992
$ws.Range("AI54").Value = 78
$ws.Range("AJ54").Value = 12446
$ws.Range("AK54").Value = 9740
$ws.Range("AL54").Value = 0
$ws.Range("AM54").Value = 1020
$ws.Range("AN54").Value = 31954
$ws.Range("AO54").Value = 0
$ws.Range("AP54").Value = 71740
$ws.Range("AQ54").Value = 17531
$ws.Range("AR54").Value = 2317
$ws.Range("AS54").Value = 136
$ws.Range("AT54").Value = 23204
$ws.Range("AU54").Value = 943
$ws.Range("AV54").Value = 0
$ws.Range("AW54").Value = 1138
$ws.Range("AX54").Value = 26472
$ws.Range("AY54").Value = 6382
$ws.Range("AZ54").Value = 149
$ws.Range("BA54").Value = 0
$ws.Range("BB54").Value = 0
$ws.Range("BC54").Value = 0
$ws.Range("BD54").Value = 30
$ws.Range("BE54").Value = 119
$ws.Range("BF54").Value = 0
$ws.Range("BG54").Value = 0
$ws.Range("BH54").Value = 0
$ws.Range("BI54").Value = 252513
